$d = $word.ActiveDocument
$quote = [char]0x201C

# --- Change 1: merge "Tombol Sebelah " / "kanan " / "bertuliskan " " into a single run ---
$txt1 = "Tombol Sebelah kanan bertuliskan " + $quote + " "
$rng1 = $d.Content
$rng1.Find.Execute($txt1, $false, $false, $false, $false, $false, $true, 1, $false, $txt1, 2) | Out-Null

# --- Change 2: merge """, Untuk " / "memajukan " / "pendorong mesin" into a single run ---
$txt2 = $quote + ", Untuk memajukan pendorong mesin"
$rng2 = $d.Content
$rng2.Find.Execute($txt2, $false, $false, $false, $false, $false, $true, 1, $false, $txt2, 2) | Out-Null

# --- Change 3: remove the existing _GoBack bookmark after "Sebelum memotong..." paragraph ---
$d.Bookmarks.Item("_GoBack").Delete()

# --- Change 4: "Ulangi langkah ke 2 dan 3 hingga..." -> "Ulangi langkah ke 2 - 4 hingga..."
#     with a new _GoBack bookmark placed right after "- 4" ---
$rng4 = $d.Content
$rng4.Find.Execute("dan 3", $false, $false, $false, $false, $false, $true, 1, $false, "- 4", 2) | Out-Null

$rngFound = $d.Content
$rngFound.Find.Execute("- 4", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Force a run split right before "- 4" with a temporary bookmark, then drop it again.
$tempRange = $d.Range($rngFound.Start, $rngFound.Start)
$d.Bookmarks.Add("zzTemp", $tempRange) | Out-Null

$realRange = $d.Range($rngFound.End, $rngFound.End)
$d.Bookmarks.Add("_GoBack", $realRange) | Out-Null

$d.Bookmarks.Item("zzTemp").Delete()
